$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 17.404364
$ws.Range("H2").Value = 52.213092
$ws.Range("I2").Value = 0.03673162149179448
$ws.Range("J2").Value = 0.03917989134808626
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.006255666666667
$ws.Range("N2").Value = 3.018767
$ws.Range("O2").Value = 0.03991953272530977
$ws.Range("P2").Value = 0.03991953272530977
$ws.Range("Q2").Value = 17.51323989972934
$ws.Range("R2").Value = 157.619159097564
$ws.Range("S2").Value = 0.001466309166195382
$ws.Range("T2").Value = 0.001564042954844011

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 17.404364
$ws.Range("H3").Value = 52.213092
$ws.Range("I3").Value = 0.03673162149179448
$ws.Range("J3").Value = 0.03917989134808626
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.646551333333332
$ws.Range("N3").Value = 22.939654
$ws.Range("O3").Value = 0.3033491053003703
$ws.Range("P3").Value = 0.3033491053003703
$ws.Range("Q3").Value = 133.0833627500187
$ws.Range("R3").Value = 1197.750264750168
$ws.Range("S3").Value = 0.01114250451576771
$ws.Range("T3").Value = 0.01188518498620769

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 17.404364
$ws.Range("H4").Value = 52.213092
$ws.Range("I4").Value = 0.03673162149179448
$ws.Range("J4").Value = 0.03917989134808626
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 16.55429333333333
$ws.Range("N4").Value = 49.66287999999999
$ws.Range("O4").Value = 0.6567313619743199
$ws.Range("P4").Value = 0.6567313619743199
$ws.Range("Q4").Value = 288.1169469361066
$ws.Range("R4").Value = 2593.05252242496
$ws.Range("S4").Value = 0.02412280780983139
$ws.Range("T4").Value = 0.02573066340703456

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 178.3379163333334
$ws.Range("H5").Value = 535.0137490000001
$ws.Range("I5").Value = 0.376379213879422
$ws.Range("J5").Value = 0.4014659877938717
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.006255666666667
$ws.Range("N5").Value = 3.018767
$ws.Range("O5").Value = 0.03991953272530977
$ws.Range("P5").Value = 0.03991953272530977
$ws.Range("Q5").Value = 179.4535388919426
$ws.Range("R5").Value = 1615.081850027483
$ws.Range("S5").Value = 0.01502488234558595
$ws.Range("T5").Value = 0.01602633463783627

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 178.3379163333334
$ws.Range("H6").Value = 535.0137490000001
$ws.Range("I6").Value = 0.376379213879422
$ws.Range("J6").Value = 0.4014659877938717
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.646551333333332
$ws.Range("N6").Value = 22.939654
$ws.Range("O6").Value = 0.3033491053003703
$ws.Range("P6").Value = 0.3033491053003703
$ws.Range("Q6").Value = 1363.670031922539
$ws.Range("R6").Value = 12273.03028730285
$ws.Range("S6").Value = 0.1141742977839794
$ws.Range("T6").Value = 0.1217843482058004

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 178.3379163333334
$ws.Range("H7").Value = 535.0137490000001
$ws.Range("I7").Value = 0.376379213879422
$ws.Range("J7").Value = 0.4014659877938717
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 16.55429333333333
$ws.Range("N7").Value = 49.66287999999999
$ws.Range("O7").Value = 0.6567313619743199
$ws.Range("P7").Value = 0.6567313619743199
$ws.Range("Q7").Value = 2952.258179437458
$ws.Range("R7").Value = 26570.32361493712
$ws.Range("S7").Value = 0.2471800337498567
$ws.Range("T7").Value = 0.263655304950235

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 86.61737833333332
$ws.Range("H8").Value = 259.852135
$ws.Range("I8").Value = 0.1828045400309692
$ws.Range("J8").Value = 0.1949889965502951
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.006255666666667
$ws.Range("N8").Value = 3.018767
$ws.Range("O8").Value = 0.03991953272530977
$ws.Range("P8").Value = 0.03991953272530977
$ws.Range("Q8").Value = 87.1592277797272
$ws.Range("R8").Value = 784.4330500175449
$ws.Range("S8").Value = 0.007297471818101475
$ws.Range("T8").Value = 0.007783869628864821

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 86.61737833333332
$ws.Range("H9").Value = 259.852135
$ws.Range("I9").Value = 0.1828045400309692
$ws.Range("J9").Value = 0.1949889965502951
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.646551333333332
$ws.Range("N9").Value = 22.939654
$ws.Range("O9").Value = 0.3033491053003703
$ws.Range("P9").Value = 0.3033491053003703
$ws.Range("Q9").Value = 662.3242297845876
$ws.Range("R9").Value = 5960.918068061289
$ws.Range("S9").Value = 0.05545359366324025
$ws.Range("T9").Value = 0.05914973764694903

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 86.61737833333332
$ws.Range("H10").Value = 259.852135
$ws.Range("I10").Value = 0.1828045400309692
$ws.Range("J10").Value = 0.1949889965502951
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 16.55429333333333
$ws.Range("N10").Value = 49.66287999999999
$ws.Range("O10").Value = 0.6567313619743199
$ws.Range("P10").Value = 0.6567313619743199
$ws.Range("Q10").Value = 1433.889488694311
$ws.Range("R10").Value = 12905.0053982488
$ws.Range("S10").Value = 0.1200534745496275
$ws.Range("T10").Value = 0.1280553892744813

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 102.6404036666667
$ws.Range("H11").Value = 307.921211
$ws.Range("I11").Value = 0.2166208691825219
$ws.Range("J11").Value = 0.231059282808824
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.006255666666667
$ws.Range("N11").Value = 3.018767
$ws.Range("O11").Value = 0.03991953272530977
$ws.Range("P11").Value = 0.03991953272530977
$ws.Range("Q11").Value = 103.2824878185374
$ws.Range("R11").Value = 929.5423903668369
$ws.Range("S11").Value = 0.008647403876316731
$ws.Range("T11").Value = 0.009223778601573453

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 102.6404036666667
$ws.Range("H12").Value = 307.921211
$ws.Range("I12").Value = 0.2166208691825219
$ws.Range("J12").Value = 0.231059282808824
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 7.646551333333332
$ws.Range("N12").Value = 22.939654
$ws.Range("O12").Value = 0.3033491053003703
$ws.Range("P12").Value = 0.3033491053003703
$ws.Range("Q12").Value = 784.8451155112214
$ws.Range("R12").Value = 7063.606039600993
$ws.Range("S12").Value = 0.06571174685590658
$ws.Range("T12").Value = 0.070091626711402

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 102.6404036666667
$ws.Range("H13").Value = 307.921211
$ws.Range("I13").Value = 0.2166208691825219
$ws.Range("J13").Value = 0.231059282808824
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 16.55429333333333
$ws.Range("N13").Value = 49.66287999999999
$ws.Range("O13").Value = 0.6567313619743199
$ws.Range("P13").Value = 0.6567313619743199
$ws.Range("Q13").Value = 1699.139350149742
$ws.Range("R13").Value = 15292.25415134768
$ws.Range("S13").Value = 0.1422617184502986
$ws.Range("T13").Value = 0.1517438774958485

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 88.82503149999999
$ws.Range("H14").Value = 177.650063
$ws.Range("I14").Value = 0.1874637554152924
$ws.Range("J14").Value = 0.1333058414989229
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.006255666666667
$ws.Range("N14").Value = 3.018767
$ws.Range("O14").Value = 0.03991953272530977
$ws.Range("P14").Value = 0.03991953272530977
$ws.Range("Q14").Value = 89.38069128872016
$ws.Range("R14").Value = 536.284147732321
$ws.Range("S14").Value = 0.007483465519110234
$ws.Range("T14").Value = 0.005321506902191209

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 88.82503149999999
$ws.Range("H15").Value = 177.650063
$ws.Range("I15").Value = 0.1874637554152924
$ws.Range("J15").Value = 0.1333058414989229
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 7.646551333333332
$ws.Range("N15").Value = 22.939654
$ws.Range("O15").Value = 0.3033491053003703
$ws.Range("P15").Value = 0.3033491053003703
$ws.Range("Q15").Value = 679.2051630497002
$ws.Range("R15").Value = 4075.230978298201
$ws.Range("S15").Value = 0.05686696248147642
$ws.Range("T15").Value = 0.04043820775001124

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 88.82503149999999
$ws.Range("H16").Value = 177.650063
$ws.Range("I16").Value = 0.1874637554152924
$ws.Range("J16").Value = 0.1333058414989229
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 16.55429333333333
$ws.Range("N16").Value = 49.66287999999999
$ws.Range("O16").Value = 0.6567313619743199
$ws.Range("P16").Value = 0.6567313619743199
$ws.Range("Q16").Value = 1470.435626793573
$ws.Range("R16").Value = 8822.613760761438
$ws.Range("S16").Value = 0.1231133274147058
$ws.Range("T16").Value = 0.08754612684672045
